# New weekly Ciboulette price record for "Vega Modelo de Temuco".
# A new row is inserted above the current row 371 (pushing the existing
# 371-411 records down to 372-412) and populated with the same values as
# the record that used to sit at row 371, except for a fresh date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("371:371").Insert()

$ws.Cells.Item(371, 1).Value = 10
$ws.Cells.Item(371, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(371, 3).Value = "La Araucanía"
$ws.Cells.Item(371, 4).Value = 45142
$ws.Cells.Item(371, 5).Value = 9
$ws.Cells.Item(371, 6).Value = 100112039
$ws.Cells.Item(371, 7).Value = "Ciboulette"
$ws.Cells.Item(371, 8).Value = "Sin especificar"
$ws.Cells.Item(371, 9).Value = "Primera"
$ws.Cells.Item(371, 10).Value = 30
$ws.Cells.Item(371, 11).Value = 7000
$ws.Cells.Item(371, 12).Value = 7000
$ws.Cells.Item(371, 13).Value = 7000
$ws.Cells.Item(371, 14).Value = "`$/docena de atados"
$ws.Cells.Item(371, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(371, 16).Value = 2333
$ws.Cells.Item(371, 17).Value = 3
$ws.Cells.Item(371, 18).Value = "Hortaliza"
